$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 08:40"

# --- Reorder country names that moved in the shared-string table ---
# These three pairs of adjacent rows swap their country-name text while the
# row's own numeric statistics stay put (rows 202/203, 208/209, 211/212).
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

$ws.Range("A211").Value = "Montserrat"
$ws.Range("A212").Value = "Seychelles"

# Rows 211/212 also carry the Seychelles/Montserrat numeric data, which swaps
# together with the names (Casos activos / Muertes columns differ).
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# --- Updated COVID-19 country statistics ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 2424492
$ws.Range("C4").Value = 324
$ws.Range("D4").Value = 1020412
$ws.Range("E4").Value = 1280604
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 123476

# India (row 7)
$ws.Range("B7").Value = 456552
$ws.Range("C7").Value = 437
$ws.Range("E7").Value = 183384

# Ucrania (row 38)
$ws.Range("B38").Value = 39014
$ws.Range("C38").Value = 940
$ws.Range("D38").Value = 17409
$ws.Range("E38").Value = 20554
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 1051

# El Salvador (row 83)
$ws.Range("D83").Value = 2847
$ws.Range("E83").Value = 2013

# Georgia (row 133)
$ws.Range("B133").Value = 914
$ws.Range("C133").Value = 3
$ws.Range("E133").Value = 132
